$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing "Text" cell type (these values look numeric/percent-like
# and would otherwise be auto-coerced to Number by Excel COM on assignment).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.930.90"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.551.25"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "304.20"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").Value = "98.15"
$ws.Range("E6").Value = "  +3.70%  "
$ws.Range("D7").Value = "0.576"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.546"
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("D10").Value = "36.89"
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("D11").Value = "0.0823"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").Value = "7.74"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").Value = "0.116"
$ws.Range("E13").Value = "  +6.19%  "
$ws.Range("D14").Value = "2.942.02"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "2.519.21"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").Value = "14.99"
$ws.Range("E16").Value = "  +5.72%  "
$ws.Range("D17").Value = "0.877"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "43.007.25"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "13.66"
$ws.Range("E19").Value = "  +3.45%  "
$ws.Range("D20").Value = "0.0₃0992"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").Value = "6.60"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").Value = "71.93"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "253.88"
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("D24").Value = "2.96"
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("E25").Value = "  -2.42%  "
$ws.Range("D26").Value = "27.92"
$ws.Range("E26").Value = "  -6.12%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "10.15"
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("D29").Value = "37.80"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("D31").Value = "6.08"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").Value = "158.46"
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "2.17"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "2.75"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("D36").Value = "19.10"
$ws.Range("E36").Value = "  +12.08%  "
$ws.Range("D37").Value = "3.31"
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("D38").Value = "26.00"
$ws.Range("E38").Value = "  +11.21%  "
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("D41").Value = "2.12"
$ws.Range("E41").Value = "  +34.04%  "
$ws.Range("D42").Value = "3.43"
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "2.090.47"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Value = "0.0305"
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "86.55"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D49").Value = "2.800.24"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "74.78"
$ws.Range("E50").Value = "  +7.47%  "
$ws.Range("D51").Value = "103.44"
$ws.Range("E51").Value = "  -1.54%  "
